$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.189.15"
$ws.Range("E2").Value = "  +2.94%  "

$ws.Range("D3").Value = "3.326.95"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.61"
$ws.Range("E5").Value = "  +3.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.74"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +2.52%  "

$ws.Range("D9").Value = "3.324.52"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.179"
$ws.Range("E10").Value = "  +1.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.42"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +3.26%  "

$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "632.69"
$ws.Range("E14").Value = "  +6.38%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.850.80"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.47"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "68.263.22"
$ws.Range("E17").Value = "  +3.08%  "

$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").Value = "3.317.35"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.70"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.903"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.64"
$ws.Range("E23").Value = "  -1.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.08"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.01"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").Value = "  +0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.58"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.47"
$ws.Range("E29").Value = "  +5.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.61"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.72"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "603.42"
$ws.Range("E32").Value = "  +7.16%  "

$ws.Range("D33").Value = "3.943.52"
$ws.Range("E33").Value = "  +3.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.98"
$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.105"
$ws.Range("E35").Value = "  +1.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.51"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.97"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  +4.49%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.129"
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  +4.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.67"

$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0688"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.37"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.339"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0416"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("E47").Value = "  +2.04%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.40"
$ws.Range("E48").Value = "  +13.72%  "

$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.97"
$ws.Range("E51").Value = "  +1.32%  "

